# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainCell($addr, $val) {
    # Plain assignment - safe for values Excel can't mistake for a number
    # (keeps the cell's existing style untouched).
    $ws.Range($addr).Value = $val
}

function Set-NumericLookingTextCell($addr, $val) {
    # The raw string would otherwise be auto-coerced to a Number by Excel's
    # type inference, so force Text formatting for the assignment, then
    # restore General formatting (the cell stays Text once it holds a
    # string value - General format on a string cell still displays/serializes
    # as text, matching the source data which has no explicit number format).
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
}

# row 2 - Bitcoin
Set-PlainCell "D2" "27.026.99"
Set-PlainCell "E2" "  -0.55%  "

# row 3 - Ethereum
Set-PlainCell "D3" "1.622.46"
Set-PlainCell "E3" "  -1.14%  "

# row 4 - TetherUSD
Set-PlainCell "E4" "  -0.19%  "

# row 5 - BNB
Set-NumericLookingTextCell "D5" "214.94"
Set-PlainCell "E5" "  -0.96%  "

# row 6 - XRP
Set-NumericLookingTextCell "D6" "0.519"
Set-PlainCell "E6" "  -1.00%  "

# row 7 - USDC
Set-PlainCell "E7" "  -0.20%  "

# row 8 / row 9 - swap Dogecoin <-> Cardano
Set-PlainCell "B8" "Cardano"
Set-PlainCell "C8" "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-NumericLookingTextCell "D8" "0.251"
Set-PlainCell "E8" "  -1.18%  "

Set-PlainCell "B9" "Dogecoin"
Set-PlainCell "C9" "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-NumericLookingTextCell "D9" "0.0628"
Set-PlainCell "E9" "  +0.22%  "

# row 10 - Solana
Set-NumericLookingTextCell "D10" "20.16"
Set-PlainCell "E10" "  +1.05%  "

# row 11 - TRON
Set-NumericLookingTextCell "D11" "0.0849"
Set-PlainCell "E11" "  +0.14%  "

# row 12 - Wrapped liquid staked Ether 2.0
Set-PlainCell "D12" "1.849.15"
Set-PlainCell "E12" "  -1.23%  "

# row 13 - Wrapped Ether
Set-PlainCell "D13" "1.618.09"
Set-PlainCell "E13" "  -1.48%  "

# row 14 - Polkadot
Set-PlainCell "E14" "  +0.13%  "

# row 15 - Polygon
Set-PlainCell "E15" "  -0.12%  "

# row 16 - Litecoin
Set-NumericLookingTextCell "D16" "64.88"
Set-PlainCell "E16" "  -3.20%  "

# row 17 - WrappedBTC
Set-PlainCell "D17" "27.000.04"
Set-PlainCell "E17" "  -0.68%  "

# row 18 - ShibaInu
Set-PlainCell "E18" "  +0.75%  "

# row 19 - BitcoinCash
Set-NumericLookingTextCell "D19" "214.23"
Set-PlainCell "E19" "  -1.72%  "

# row 20 - Dai
Set-PlainCell "E20" "  -0.12%  "

# row 21 - Chainlink
Set-NumericLookingTextCell "D21" "6.87"
Set-PlainCell "E21" "  -1.12%  "

# row 22 - Uniswap
Set-PlainCell "E22" "  -1.20%  "

# row 23 - Toncoin
Set-PlainCell "E23" "  -5.23%  "

# row 24 - Avalanche
Set-PlainCell "E24" "  -0.42%  "

# row 25 - Monero
Set-NumericLookingTextCell "D25" "148.28"
Set-PlainCell "E25" "  +0.72%  "

# row 26 - BinanceUSD
Set-PlainCell "E26" "  -0.20%  "

# row 27 - Cosmos
Set-PlainCell "E27" "  -0.80%  "

# row 28 - Stellar
Set-PlainCell "E28" "  -2.42%  "

# row 29 - EthereumClassic
Set-NumericLookingTextCell "D29" "15.58"
Set-PlainCell "E29" "  -0.78%  "

# row 30 - Hedera
Set-NumericLookingTextCell "D30" "0.0513"
Set-PlainCell "E30" "  +1.06%  "

# row 31 - PancakeSwap
Set-PlainCell "E31" "  -0.87%  "

# row 32 - ImmutableX
Set-NumericLookingTextCell "D32" "0.755"
Set-PlainCell "E32" "  +37.70%  "

# row 33 - Filecoin
Set-PlainCell "E33" "  -0.80%  "

# row 34 - InternetComputer(DFINITY)
Set-PlainCell "E34" "  -0.23%  "

# row 35 - Maker
Set-PlainCell "D35" "1.349.10"
Set-PlainCell "E35" "  +3.80%  "

# row 36 - LidoDAOToken
Set-PlainCell "E36" "  -0.34%  "

# row 37 - HuobiToken
Set-PlainCell "E37" "  -0.38%  "

# row 38 - VeChain
Set-PlainCell "E38" "  +0.66%  "

# row 39 - ARBITRUM
Set-NumericLookingTextCell "D39" "0.847"
Set-PlainCell "E39" "  -1.06%  "

# row 40 - PaxDollar
Set-PlainCell "E40" "  -0.17%  "

# row 41 - TrustWalletToken
Set-NumericLookingTextCell "D41" "0.800"
Set-PlainCell "E41" "  -1.19%  "

# row 42 - MXToken
Set-PlainCell "E42" "  -0.19%  "

# row 43 - Aave
Set-NumericLookingTextCell "D43" "65.12"
Set-PlainCell "E43" "  +5.52%  "

# row 44 - FraxShare
Set-PlainCell "E44" "  +0.44%  "

# row 45 - RocketPoolETH
Set-PlainCell "D45" "1.760.64"

# row 46 - WEMIXToken
Set-NumericLookingTextCell "D46" "0.882"
Set-PlainCell "E46" "  +31.77%  "

# row 47 - Quant
Set-NumericLookingTextCell "D47" "90.14"
Set-PlainCell "E47" "  -1.79%  "

# row 48 - RenderToken
Set-PlainCell "E48" "  +2.51%  "

# row 49 - BabyDogeCoin
Set-PlainCell "E49" "  -0.89%  "

# row 50 - Algorand
Set-NumericLookingTextCell "D50" "0.101"
Set-PlainCell "E50" "  +5.41%  "

# row 51 - Cronos
Set-PlainCell "E51" "  +0.38%  "
